# Refresh the crypto price/volume table with the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" cells hold digit-grouped text (e.g. "542.99") that Excel
# would otherwise silently reinterpret as a number when assigned via .Value.
# Flip those cells to text format first so the literal string sticks, then
# clear the formatting delta back off so no stray number style lingers.
$textCells = @("D5","D6","D8","D12","D15","D19","D20","D21","D23","D24","D25","D27","D28","D31","D33","D34","D35","D36","D37","D39","D41","D42","D43","D44","D46","D47","D48","D49")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "59.180.87"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.510.86"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("D5").Value = "542.99"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "144.18"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").Value = "2.536.99"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "5.55"
$ws.Range("E12").Value = "  +4.10%  "
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "2.954.54"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").Value = "23.70"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "59.101.07"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "2.527.31"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "4.30"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").Value = "325.34"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").Value = "5.81"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("D24").Value = "61.86"
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("D25").Value = "0.439"
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "7.98"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "6.66"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  -6.18%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "1.48"
$ws.Range("E33").Value = "  +8.34%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "157.99"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "18.69"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "4.38"
$ws.Range("E37").Value = "  -2.98%  "
$ws.Range("E38").Value = "  -5.99%  "
$ws.Range("D39").Value = "5.65"
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").Value = "299.15"
$ws.Range("E41").Value = "  -6.18%  "
$ws.Range("D42").Value = "3.71"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "0.822"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("D46").Value = "10.78"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").Value = "0.0930"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").Value = "18.78"
$ws.Range("E48").Value = "  +2.11%  "
$ws.Range("D49").Value = "123.11"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("E51").Value = "  -0.35%  "

foreach ($addr in $textCells) { $ws.Range($addr).ClearFormats() }
